# MT_RR_10.xlsx edit script
# Adds an "Elapsed Time" column (based on End_time) and splits the single
# "Fairness" column into three: Fairness(QWT), Fairness(RT), Fairness(ET).
# Also nudges the default font from Calibri/Cambria to Arial/Times New Roman.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new columns -------------------------------------------
# Original layout:  ... H(Response Time) I(Turnaround Time) J(Throughput)
#                    K(CPU Usage) L(Fairness) M(Context Switch)
# Target layout:     ... H(Response Time) I(Elapsed Time) J(Turnaround Time)
#                    K(Throughput) L(CPU Usage) M(Fairness(QWT))
#                    N(Fairness(RT)) O(Fairness(ET)) P(Context Switch)

# New column for "Elapsed Time", right after Response Time.
$ws.Columns("I:I").Insert()

# New column right before the existing "Fairness" column (now shifted to M).
$ws.Columns("M:M").Insert()

# New column right after the (now shifted) "Fairness" column (now at N).
$ws.Columns("O:O").Insert()

# --- 2. Headers ------------------------------------------------------------
# (order matters for shared-string allocation order, matching the original
# authoring sequence: the three Fairness headers were introduced before the
# Elapsed Time header)
$ws.Range("M1").Value = "Fairness(QWT)"
$ws.Range("N1").Value = "Fairness(RT)"
$ws.Range("O1").Value = "Fairness(ET)"
$ws.Range("I1").Value = "Elapsed Time"

# --- 3. Formulas -------------------------------------------------------------
# Elapsed Time per row = (End_time - min(Start_time)) / 100
$ws.Range("I2:I11").Formula = "=(C2-`$B`$12)/100"

# Fairness(QWT): stddev of Wait Time (column G)
$ws.Range("M12").Formula = "=_xlfn.STDEV.P(G2:G11)"

# Fairness(RT): stddev of Real_time (column D) -- this is the original
# "Fairness" formula, kept as-is (now living in column N).

# Fairness(ET): stddev of Elapsed Time (column I)
$ws.Range("O12").Formula = "=_xlfn.STDEV.P(I2:I11)"

# --- 4. Cosmetic: default font (Calibri/Cambria -> Arial/Times New Roman) --
$ws.Cells.Font.Name = "Arial"

# --- 5. View state: selection moved to O14 ----------------------------------
$ws.Range("O14").Select()
